# This script reproduces a weekly data refresh of the "Fruta, Comercializadora
# del Agro de Limarí - Arándano (blue)" sheet: for most data rows the Fecha
# (D), Volumen (M), Precio mínimo (N), Precio máximo (O), Precio promedio
# ponderado (P) and Precio $/Kg (S) values were replaced with another row's
# values (a reshuffle of daily records), while rows 3 and 6 stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @{
    2  = @{ D = 44818; M = 200;  N = 11000; O = 12000; P = 11500; S = 5750 }
    4  = @{ D = 45203; M = 1000; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    5  = @{ D = 44482; M = 240;  N = 10000; O = 11000; P = 10500; S = 5250 }
    7  = @{ D = 44874; M = 300;  N = 7500;  O = 8000;  P = 7750;  S = 3875 }
    8  = @{ D = 45211; M = 200;  N = 10000; O = 11000; P = 10500; S = 5250 }
    9  = @{ D = 44882; M = 440;  N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    10 = @{ D = 45204; M = 400;  N = 9000;  O = 10000; P = 9500;  S = 4750 }
    11 = @{ D = 44490; M = 400;  N = 9500;  O = 10000; P = 9750;  S = 4875 }
    12 = @{ D = 44517; M = 400;  N = 5500;  O = 6000;  P = 5750;  S = 2875 }
    13 = @{ D = 44489; M = 160;  N = 9500;  O = 10000; P = 9750;  S = 4875 }
    14 = @{ D = 44895; M = 240;  N = 3000;  O = 3500;  P = 3250;  S = 1625 }
    15 = @{ D = 44881; M = 440;  N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    16 = @{ D = 44475; M = 240;  N = 11000; O = 12000; P = 11500; S = 5750 }
    17 = @{ D = 44454; M = 160;  N = 12000; O = 13000; P = 12500; S = 6250 }
    18 = @{ D = 44819; M = 240;  N = 11000; O = 12000; P = 11500; S = 5750 }
    19 = @{ D = 44889; M = 460;  N = 3500;  O = 4000;  P = 3750;  S = 1875 }
    20 = @{ D = 44455; M = 200;  N = 12000; O = 13000; P = 12500; S = 6250 }
    21 = @{ D = 44461; M = 200;  N = 11000; O = 12000; P = 11500; S = 5750 }
}

foreach ($row in $targets.Keys) {
    $vals = $targets[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N - Precio mínimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Precio máximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals.S   # S - Precio $/Kg
}
